$d = $word.ActiveDocument

$pairs = @(
    @{old = "442×4="; new = "445×4="},
    @{old = "959×9="; new = "395×6="},
    @{old = "240×6="; new = "239×3="},
    @{old = "476×5="; new = "942×5="},
    @{old = "303×8="; new = "433×7="},
    @{old = "623×9="; new = "203×6="},
    @{old = "320×8="; new = "681×8="},
    @{old = "606×2="; new = "248×7="},
    @{old = "374×7="; new = "998×7="},
    @{old = "638×7="; new = "802×6="},
    @{old = "709×9="; new = "776×9="},
    @{old = "375×2="; new = "234×6="},
    @{old = "737×4="; new = "177×4="},
    @{old = "468×9="; new = "796×5="},
    @{old = "233×4="; new = "990×5="},
    @{old = "548×4="; new = "249×5="},
    @{old = "606×6="; new = "120×8="},
    @{old = "643×7="; new = "176×7="},
    @{old = "367×9="; new = "867×8="},
    @{old = "797×7="; new = "331×4="},
    @{old = "632×8="; new = "851×9="},
    @{old = "998×2="; new = "604×7="},
    @{old = "406×9="; new = "224×4="},
    @{old = "874×9="; new = "494×4="},
    @{old = "310×6="; new = "460×4="}
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.new, 2) | Out-Null
}
